# Update the forecast error table values (B2:F10) with corrected figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -0.04774544758327767
$ws.Range("C2").Value = 0.4624696470792148
$ws.Range("D2").Value = 0.8033151732794445
$ws.Range("E2").Value = 0.8962785132309289
$ws.Range("F2").Value = 0.9287915074979147

$ws.Range("B3").Value = -0.04718447073015461
$ws.Range("C3").Value = 0.5635575896698916
$ws.Range("D3").Value = 0.8690252791348545
$ws.Range("E3").Value = 0.9322152536484557
$ws.Range("F3").Value = 0.9690367108258421

$ws.Range("B4").Value = -0.03314431376102876
$ws.Range("C4").Value = 0.5287269987382343
$ws.Range("D4").Value = 0.6606252920084058
$ws.Range("E4").Value = 0.8127885899841396
$ws.Range("F4").Value = 0.8482238639804189

$ws.Range("B5").Value = -0.1426010309797199
$ws.Range("C5").Value = 0.5484211637281492
$ws.Range("D5").Value = 0.6964982596166069
$ws.Range("E5").Value = 0.8345647126596037
$ws.Range("F5").Value = 0.862426533762813

$ws.Range("B6").Value = -0.2189447014214337
$ws.Range("C6").Value = 0.5233672663903658
$ws.Range("D6").Value = 0.6294409302107248
$ws.Range("E6").Value = 0.7933731342884789
$ws.Range("F6").Value = 0.8038132369664293

$ws.Range("B7").Value = 0.05889770248139205
$ws.Range("C7").Value = 0.3795808076999097
$ws.Range("D7").Value = 0.2585567114336377
$ws.Range("E7").Value = 0.5084847209441379
$ws.Range("F7").Value = 0.5356993033274895

$ws.Range("B8").Value = -0.01074580271062597
$ws.Range("C8").Value = 0.4136063680922605
$ws.Range("D8").Value = 0.2412788381939795
$ws.Range("E8").Value = 0.4912014232409954
$ws.Range("F8").Value = 0.5379554248278389

$ws.Range("B9").Value = 0.1748454189159413
$ws.Range("C9").Value = 0.4191804083912449
$ws.Range("D9").Value = 0.2366906289677956
$ws.Range("E9").Value = 0.4865086114014793
$ws.Range("F9").Value = 0.556039173689999

$ws.Range("B10").Value = 0.7444949522781628
$ws.Range("C10").Value = 0.7444949522781628
$ws.Range("D10").Value = 0.5542727339676639
$ws.Range("E10").Value = 0.7444949522781628
